# Update the MSME Mozambique summary figures to their more precise values.
# These cells hold numeric-looking figures but are stored as TEXT
# (e.g. "1.1" -> "1.13"), so we force a Text number format before writing
# the new value to stop Excel from auto-converting the string to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B13" = "1.13"
    "C13" = "0.26"
    "D13" = "1.39"
    "B14" = "14.13"
    "C14" = "28.78"
    "D14" = "42.91"
    "B16" = "71.13"
    "C16" = "16.65"
    "D16" = "87.77"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
